# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the
# zh-cn and de-de handback-status report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 11:21:46"
$wsZhCn.Range("H2").Value = "2016-03-17 11:22:04"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 11:21:49"
$wsDeDe.Range("H2").Value = "2016-03-17 11:22:09"
